$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: Fix "Issac Schultz" -> "Isaac Schultz" in the second
# paragraph, reproducing the exact run split from the diff:
#   <w:r>Is</w:r><w:r>a</w:r><w:bookmarkStart .../><w:bookmarkEnd .../><w:r>ac Schultz</w:r>
# and drop the old spell-check proofErr markers around the old
# "Issac" run.
# ------------------------------------------------------------------

$nameRange = $d.Paragraphs(2).Range
$nameStart = $nameRange.Start
$nameEnd = $nameRange.End

# Delete the whole paragraph's content *including* its end-of-paragraph
# mark; this properly clears away the old <w:proofErr/> markers that a
# plain text replacement would otherwise leave behind as orphans.
$d.Range($nameStart, $nameEnd).Delete()

# Re-create an (empty) paragraph in the same spot, inheriting the
# paragraph formatting of the paragraph that follows it (same as the
# original "Issac Schultz" paragraph had).
$d.Paragraphs(2).Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(2).Range
$insStart = $newPara.Start

# Insert the corrected, fully spelled out name as a single run first.
$d.Range($insStart, $insStart).Text = "Isaac Schultz"

# Move the document's "_GoBack" bookmark (it currently sits inside the
# long paragraph below) to sit between "Isa" and "ac Schultz" -- this
# also splits that text into two runs ("Isa" | "ac Schultz") and
# removes the bookmark from its old location since bookmark names are
# unique.
$d.Bookmarks.Add("_GoBack", $d.Range($insStart + 3, $insStart + 3))

# Force a further run split of "Isa" into "Is" | "a" by briefly adding
# a temporary bookmark between them, then removing it again -- the run
# break that the bookmark forces stays in place even after the
# bookmark itself is deleted.
$d.Bookmarks.Add("zzTmpSplit", $d.Range($insStart + 2, $insStart + 2))
$d.Bookmarks.Item("zzTmpSplit").Delete()

# ------------------------------------------------------------------
# Part 2: Undo the old run split in the long paragraph that the
# "_GoBack" bookmark used to create (".. we creat" | "ed an abstract
# parent class ..") now that the bookmark has moved away from there;
# merge the text back into a single run.
# ------------------------------------------------------------------

$fullText = $d.Content.Text
$splitIdx = $fullText.IndexOf("ed an abstract parent class")
$mergeRange = $d.Range($splitIdx - 8, $splitIdx + "ed an abstract parent class".Length)
$mergeRange.Text = "we created an abstract parent class"
